# Changelog update: add version 1.3.2 entry (commit: "adapt changelog to version 1.3.2")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changeText = "Improvements:`n- Change line style of absolute prices on coinprice evaluation`n- Add MyDefichain as another category to masternode evaluation`n- Reduction digits of hover label on DefiChain Promo evaluation"

# New row 15: Date | Version | Changes
$ws.Range("A15").Value = 44340
$ws.Range("B15").Value = "1.3.2"
$ws.Range("C15").Value = $changeText

# Match formatting used by the other changelog rows (column default styles
# already carry the date / version number formats, so only wrap text needs
# to be applied explicitly to the Changes column).
$ws.Range("C15").WrapText = $true

# Row height matches the 4 wrapped lines of the new changelog entry (15pt/line)
$ws.Rows.Item(15).RowHeight = 60

# Move the active selection down to the next empty row, as Excel would leave it
# after the new row was entered.
[void]$ws.Range("C16").Select()
